# Apply the commit's change: "remove Gamelogic project, modify SLG building config"
# Concretely, the worksheet gains a second data column (B) describing, for each
# effect-id row in column A, which "setting" resource it belongs to:
#   - row 1 (header) already has B1 = "Atlas_ResID"; unchanged
#   - row 2 (EFT_INFO)      -> B2  = "msg_icon"
#   - rows 3-15 (the other EFT_* effects) -> B = "Ssetting"
# Column A itself is left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "msg_icon"

$ssettingRows = 3..15
foreach ($r in $ssettingRows) {
    $ws.Range("B$r").Value = "Ssetting"
}

# Match the final selection left behind in the saved workbook.
$ws.Range("E14").Select() | Out-Null
